# Apply the edit described by the diff:
#  - Metadata!B8 (Date value) changes from 2025-07-16T13:52:06+00:00 to 2025-07-17T14:35:50+00:00
#  - Metadata!B5 (Title value) gets filled in with the same text used for Name (B4): "SurspecialiteTransversale"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"

# Fill in the Title value (row 5, column B) reusing the Name value text
$ws.Range("B5").Value = "SurspecialiteTransversale"
